{"js": "// Remove the leftover \"_GoBack\" bookmark.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Locate the \">>>  your stuff after this line >>>\" paragraph (it is\n// currently split across several runs with proofErr markers) and\n// normalize it to a single run, then add the new paragraph about\n// version control right after it.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"your stuff after this line\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target) {\n  target.clear();\n  target.insertText(\">>>  your stuff after this line >>>\", Word.InsertLocation.start);\n  target.insertParagraph(\n    \"Version control is a way to keep a track of the changes in the code so that if something goes wrong, we can make comparisons in different code versions and revert to any previous version that we want. It is very much required where multiple developers are continuously working on /changing the source code.\",\n    Word.InsertLocation.after\n  );\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove the leftover \"_GoBack\" bookmark.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Find the paragraph that holds \">>>  your stuff after this line >>>\"\n# (currently split across several runs with grammar-check proofErr\n# markers) and normalize it down to a single run with identical text.\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*your stuff after this line*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $full = \">>>  your stuff after this line >>>\"\n\n    # Delete the paragraph's text (but keep its paragraph mark) then\n    # re-insert the text as one clean run.\n    $body = $d.Range($target.Range.Start, $target.Range.End - 1)\n    $body.Delete()\n\n    $target = $d.Paragraphs.Item($i)\n    $insertPoint = $d.Range($target.Range.Start, $target.Range.Start)\n    $insertPoint.InsertAfter($full)\n\n    # Add the new paragraph about version control right after it.\n    $target = $d.Paragraphs.Item($i)\n    $target.Range.InsertParagraphAfter()\n\n    $newPara = $d.Paragraphs.Item($i + 1)\n    $newInsertPoint = $d.Range($newPara.Range.Start, $newPara.Range.Start)\n    $newInsertPoint.InsertAfter(\"Version control is a way to keep a track of the changes in the code so that if something goes wrong, we can make comparisons in different code versions and revert to any previous version that we want. It is very much required where multiple developers are continuously working on /changing the source code.\")\n}\n"}
